# Added scoreboard auto ordering functionality
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Pink",   "test", 6,  "2024-11-08 04:24:11"),
    @("Blue",   "Ace",  6,  "2024-11-08 04:24:30"),
    @("Purple", "nice", 8,  "2024-11-08 04:24:54"),
    @("Blue",   "nice", -6, "2024-11-08 04:25:08")
)

$row = 4
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    # Column B (Individual Name) is left blank for these auto-logged rows
    $ws.Cells.Item($row, 3).Value = $entry[1]
    $ws.Cells.Item($row, 4).Value = $entry[2]
    $ws.Cells.Item($row, 5).Value = $entry[3]
    $row++
}
